$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - update "想去人数" (F column) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 282
$ws1.Range("F6").Value = 12
$ws1.Range("F7").Value = 298
$ws1.Range("F8").Value = 8031
$ws1.Range("F9").Value = 73
$ws1.Range("F12").Value = 107
$ws1.Range("F15").Value = 20
$ws1.Range("F19").Value = 694
$ws1.Range("F20").Value = 24
$ws1.Range("F21").Value = 77

# Sheet "演出" (Performance) - update "想去人数" (F column) value
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 9

# Sheet "全部类型" (All types) - update "想去人数" (F column) values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 282
$ws4.Range("F6").Value = 12
$ws4.Range("F7").Value = 298
$ws4.Range("F8").Value = 8031
$ws4.Range("F9").Value = 73
$ws4.Range("F12").Value = 107
$ws4.Range("F15").Value = 20
$ws4.Range("F19").Value = 694
$ws4.Range("F20").Value = 24
$ws4.Range("F21").Value = 77
$ws4.Range("F22").Value = 9

$wb.Save()
